$wb = $excel.ActiveWorkbook

# 1. Rename header in "Weekly Quantity" sheet (B1)
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2. Rename header in "Monthly Trend" sheet (B1)
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add a new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins used by the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# 4. Write header row
$header = New-Object "object[,]" 1,4
$header[0,0] = "ds"
$header[0,1] = "PO_Forecast"
$header[0,2] = "yhat_lower"
$header[0,3] = "yhat_upper"
$wsForecast.Range("A1:D1").Value = $header

# 5. Write data rows 2..55
$data = New-Object "object[,]" 54,4
$data[0,0] = 45060.99999999999
$data[0,1] = 40
$data[0,2] = -31.6307078657743
$data[0,3] = 104.8908076813847
$data[1,0] = 45067.99999999999
$data[1,1] = 40
$data[1,2] = -27.53200517609202
$data[1,3] = 109.1763182659417
$data[2,0] = 45074.99999999999
$data[2,1] = 40
$data[2,2] = -27.83675379159009
$data[2,3] = 107.2387939751716
$data[3,0] = 45081.99999999999
$data[3,1] = 40
$data[3,2] = -28.01702101728591
$data[3,3] = 110.9160831714085
$data[4,0] = 45088.99999999999
$data[4,1] = 41
$data[4,2] = -23.71912317669854
$data[4,3] = 109.5331256565863
$data[5,0] = 45095.99999999999
$data[5,1] = 41
$data[5,2] = -23.82315247235982
$data[5,3] = 111.2422423428962
$data[6,0] = 45102.99999999999
$data[6,1] = 41
$data[6,2] = -26.40159494667416
$data[6,3] = 109.3890377466603
$data[7,0] = 45109.99999999999
$data[7,1] = 41
$data[7,2] = -28.02056293934662
$data[7,3] = 106.9661910603202
$data[8,0] = 45116.99999999999
$data[8,1] = 42
$data[8,2] = -25.38562177186942
$data[8,3] = 111.2873907078682
$data[9,0] = 45137.99999999999
$data[9,1] = 42
$data[9,2] = -26.79562148563044
$data[9,3] = 104.2504740051727
$data[10,0] = 45186.99999999999
$data[10,1] = 44
$data[10,2] = -26.63836822762855
$data[10,3] = 111.8232931581874
$data[11,0] = 45207.99999999999
$data[11,1] = 45
$data[11,2] = -19.17052986516014
$data[11,3] = 111.3450184967649
$data[12,0] = 45214.99999999999
$data[12,1] = 45
$data[12,2] = -23.3900389608616
$data[12,3] = 112.8589671933831
$data[13,0] = 45221.99999999999
$data[13,1] = 45
$data[13,2] = -22.14933929441103
$data[13,3] = 107.3707306468695
$data[14,0] = 45228.99999999999
$data[14,1] = 45
$data[14,2] = -25.1210130608463
$data[14,3] = 108.3978529485566
$data[15,0] = 45242.99999999999
$data[15,1] = 46
$data[15,2] = -22.52706699764084
$data[15,3] = 114.309330293307
$data[16,0] = 45249.99999999999
$data[16,1] = 46
$data[16,2] = -20.78435061896659
$data[16,3] = 114.3346438035721
$data[17,0] = 45256.99999999999
$data[17,1] = 46
$data[17,2] = -16.17492165631649
$data[17,3] = 114.0327601253962
$data[18,0] = 45270.99999999999
$data[18,1] = 47
$data[18,2] = -17.1798324946428
$data[18,3] = 112.021983848564
$data[19,0] = 45277.99999999999
$data[19,1] = 47
$data[19,2] = -22.63399391621601
$data[19,3] = 110.6769202831366
$data[20,0] = 45298.99999999999
$data[20,1] = 48
$data[20,2] = -20.02147788821997
$data[20,3] = 114.1442030371193
$data[21,0] = 45305.99999999999
$data[21,1] = 48
$data[21,2] = -21.48287389853882
$data[21,3] = 115.0343390358398
$data[22,0] = 45312.99999999999
$data[22,1] = 48
$data[22,2] = -15.16310028943466
$data[22,3] = 115.3843580693732
$data[23,0] = 45326.99999999999
$data[23,1] = 49
$data[23,2] = -19.85347674603991
$data[23,3] = 112.8557616841222
$data[24,0] = 45333.99999999999
$data[24,1] = 49
$data[24,2] = -21.19252856000998
$data[24,3] = 111.3570131814147
$data[25,0] = 45340.99999999999
$data[25,1] = 49
$data[25,2] = -22.75134050406884
$data[25,3] = 113.2752754165414
$data[26,0] = 45347.99999999999
$data[26,1] = 49
$data[26,2] = -16.66219135880507
$data[26,3] = 120.0544461056972
$data[27,0] = 45354.99999999999
$data[27,1] = 50
$data[27,2] = -21.18636076558043
$data[27,3] = 113.2394039712714
$data[28,0] = 45361.99999999999
$data[28,1] = 50
$data[28,2] = -19.48888187425717
$data[28,3] = 114.2249445341996
$data[29,0] = 45368.99999999999
$data[29,1] = 50
$data[29,2] = -16.48924792807441
$data[29,3] = 118.9555351564276
$data[30,0] = 45375.99999999999
$data[30,1] = 50
$data[30,2] = -17.27787944615373
$data[30,3] = 118.3174808266847
$data[31,0] = 45382.99999999999
$data[31,1] = 51
$data[31,2] = -19.44981881659761
$data[31,3] = 116.5666991282313
$data[32,0] = 45389.99999999999
$data[32,1] = 51
$data[32,2] = -15.68021972419508
$data[32,3] = 114.2928210741459
$data[33,0] = 45396.99999999999
$data[33,1] = 51
$data[33,2] = -18.0308421315093
$data[33,3] = 118.095954465014
$data[34,0] = 45410.99999999999
$data[34,1] = 52
$data[34,2] = -20.26477465523827
$data[34,3] = 119.8622007493068
$data[35,0] = 45417.99999999999
$data[35,1] = 52
$data[35,2] = -17.09044178576201
$data[35,3] = 120.4151405881365
$data[36,0] = 45424.99999999999
$data[36,1] = 52
$data[36,2] = -20.3958589585968
$data[36,3] = 115.8493988069367
$data[37,0] = 45459.99999999999
$data[37,1] = 53
$data[37,2] = -13.16491665693214
$data[37,3] = 115.9959379090654
$data[38,0] = 45522.99999999999
$data[38,1] = 55
$data[38,2] = -9.524910999788826
$data[38,3] = 120.5763475166941
$data[39,0] = 45550.99999999999
$data[39,1] = 56
$data[39,2] = -12.83393130374947
$data[39,3] = 124.9427837682454
$data[40,0] = 45557.99999999999
$data[40,1] = 57
$data[40,2] = -13.62461734736682
$data[40,3] = 125.1116041722581
$data[41,0] = 45564.99999999999
$data[41,1] = 57
$data[41,2] = -7.600480929066136
$data[41,3] = 124.5917216797524
$data[42,0] = 45571.99999999999
$data[42,1] = 57
$data[42,2] = -13.17821025120075
$data[42,3] = 121.9966608913311
$data[43,0] = 45578.99999999999
$data[43,1] = 57
$data[43,2] = -10.75313675296103
$data[43,3] = 125.1363124183471
$data[44,0] = 45585.99999999999
$data[44,1] = 57
$data[44,2] = -6.410320510393823
$data[44,3] = 123.8581620516149
$data[45,0] = 45599.99999999999
$data[45,1] = 58
$data[45,2] = -4.471589852048992
$data[45,3] = 129.7372753970445
$data[46,0] = 45606.99999999999
$data[46,1] = 58
$data[46,2] = -11.41635087608994
$data[46,3] = 127.6861734977909
$data[47,0] = 45613.99999999999
$data[47,1] = 58
$data[47,2] = -8.515205219869808
$data[47,3] = 124.8538464043907
$data[48,0] = 45620.99999999999
$data[48,1] = 59
$data[48,2] = -12.63345069058207
$data[48,3] = 123.0805460901917
$data[49,0] = 45627.99999999999
$data[49,1] = 59
$data[49,2] = -15.20740224643026
$data[49,3] = 125.6013433772248
$data[50,0] = 45634.99999999999
$data[50,1] = 59
$data[50,2] = -8.936896782698804
$data[50,3] = 126.3802796119275
$data[51,0] = 45641.99999999999
$data[51,1] = 59
$data[51,2] = -3.487277976154745
$data[51,3] = 126.7550357560466
$data[52,0] = 45648.99999999999
$data[52,1] = 60
$data[52,2] = -5.282236966446
$data[52,3] = 130.8925269734979
$data[53,0] = 45655.99999999999
$data[53,1] = 60
$data[53,2] = -3.499108669134637
$data[53,3] = 130.7208411393327
$wsForecast.Range("A2:D55").Value = $data

# 6. Apply header formatting (bold, border, centered) matching existing sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# 7. Apply date-time number format to column A (matching existing date columns)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A55").PasteSpecial(-4122)

$excel.CutCopyMode = $false
